# Product List update (ProductsData sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductsData")

# Row 2: GPI-0181 (Epson ink) -> GPT-0516 (Konica Minolta toner)
$ws.Range("A2").Value = "GPT-0516"
$ws.Range("B2").Value = "KONICA MINOLTA TONER 2300 BLACK ΣΥΜΒΑΤΟ 4000 ΣΕΛΙΔΕΣ"
$ws.Range("C2").Value = 3

# Row 3: GPI-0017 (Epson ink) -> GPI-0079 (HP ink); quantity kept as text "2"
$ws.Range("A3").Value = "GPI-0079"
$ws.Range("B3").Value = "HP INK No 351XL - CB338EE COLOR ΣΥΜΒΑΤΟ 18ml"
$ws.Range("C3").Value = "'2"

# Rows 4 and 5 (GPT-0219 / GPT-0380) are removed entirely, shrinking the
# used range down to A1:C3.
$ws.Range("A4:A5").EntireRow.Delete()
